$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The perk table (A1:D57) is re-sorted: primary key column A (Ability set),
# secondary key column C (UIPerk_* icon name), both ascending. This mirrors
# the <sortState>/<sortCondition> bookkeeping that Excel records after a
# Data > Sort operation on this range.
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("A1:A57")) | Out-Null
$sortObj.SortFields.Add($ws.Range("C1:C57")) | Out-Null
$sortObj.SetRange($ws.Range("A1:D57"))
$sortObj.Header = -4142
$sortObj.Apply()

# Fill in newly-reviewed icon statuses in column D for several Engineer &
# Dragoon perks (most of the commit's actual content change).
$ws.Range("D2").Value = "As-is"
$ws.Range("D3").Value = "OK"
$ws.Range("D6").Value = "OK"
$ws.Range("D8").Value = "OK"
$ws.Range("D12").Value = "As-is"
$ws.Range("D13").Value = "OK"
$ws.Range("D14").Value = "OK"
$ws.Range("D15").Value = "OK"
$ws.Range("D16").Value = "OK"
$ws.Range("D17").Value = "OK"
$ws.Range("D18").Value = "OK"
$ws.Range("D19").Value = "OK"
$ws.Range("D20").Value = "OK"
$ws.Range("D22").Value = "OK"
$ws.Range("D23").Value = "Check"
$ws.Range("D25").Value = "OK"
$ws.Range("D27").Value = "Check"

# Move the active selection to D20, matching where the edits were being made.
$ws.Range("D20").Select()
